$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 126: staining
$ws.Cells.Item(126, 1).Value = "http://purl.obolibrary.org/obo/OBI_0302887"
$ws.Cells.Item(126, 2).Value = "staining"
$ws.Cells.Item(126, 3).Value = "y"
$ws.Cells.Item(126, 5).Value = "yes"

# Row 127: IHC-stained fixed tissue slide specimen
$ws.Cells.Item(127, 1).Value = "http://purl.obolibrary.org/obo/OBI_0002126"
$ws.Cells.Item(127, 2).Value = "IHC-stained fixed tissue slide specimen"
$ws.Cells.Item(127, 3).Value = "y"

# Row 128: H&E-stained fixed tissue slide specimen
$ws.Cells.Item(128, 1).Value = "http://purl.obolibrary.org/obo/OBI_0002125"
$ws.Cells.Item(128, 2).Value = "H&E-stained fixed tissue slide specimen"
$ws.Cells.Item(128, 3).Value = "y"

# Update the view: scroll/select near the newly added rows and zoom to 150%
$ws.Range("C128").Select()
$excel.ActiveWindow.Zoom = 150
